$d = $word.ActiveDocument

# The document currently has a single paragraph:
#   "DỰ ÁN PHÂN TÍCH DỮ LIỆU"
# immediately followed (in the same paragraph) by the "_GoBack" bookmark.
# We need to split this into two paragraphs, keeping the title paragraph as-is
# and adding a new, centered & bold paragraph "Họ và tên: Nguyễn Phước" that
# ends up owning the "_GoBack" bookmark (matching the target OOXML).

# Remove the existing "_GoBack" bookmark first; we'll recreate it (via raw XML)
# at the correct spot after the new paragraph is created, since this engine
# does not relocate zero-width bookmarks across paragraph splits automatically.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the end of the title paragraph's text (just before its paragraph mark).
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$insertionPoint = $d.Range($titleRange.Start, $titleRange.End - 1)
$insertionPoint.Collapse(0)  # wdCollapseEnd

# Insert a paragraph break, then the new bold, centered line of text.
$insertionPoint.Font.Bold = 1
$insertionPoint.InsertAfter([char]13 + "Họ và tên: Nguyễn Phước")

# Replace the newly created second paragraph's contents with an equivalent
# OOXML fragment that also carries the "_GoBack" bookmark, so the bookmark
# ends up inside this paragraph (right after the text, before the paragraph
# mark) instead of staying behind in the title paragraph.
$newPara = $d.Paragraphs(2)
$newParaRange = $newPara.Range
$newParaTextRange = $d.Range($newParaRange.Start, $newParaRange.End - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Họ và tên: Nguyễn Phước</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newParaTextRange.InsertXML($xml)
